$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix E8 value (0.091 -> 0.09)
$ws.Range("E8").Value = 0.09

# Row 10
$ws.Range("A10").Value = "voor oplossing checken voordat bord aan queue wordt toegevoegd"
$ws.Range("C10").Value = 0.83
$ws.Range("D10").Value = 0.14
$ws.Range("E10").Value = 0.07
$ws.Range("F10").Value = 24.78
$ws.Range("F10").NumberFormat = "0.0"

# Row 11
$ws.Range("A11").Value = "vermijden occupiedBy()"
$ws.Range("C11").Value = 0.8
$ws.Range("C11").NumberFormat = "0.00"
$ws.Range("D11").Value = 0.13
$ws.Range("E11").Value = 0.07
$ws.Range("F11").Value = 25.2
$ws.Range("F11").NumberFormat = "0.0"

# Row 12 - only F12 carries formatting, no values
$ws.Range("F12").NumberFormat = "0.0"

# Row 13
$ws.Range("A13").Value = "Heapq: heuristic = numMoves"
$ws.Range("C13").Value = 0.89
$ws.Range("D13").Value = 0.14
$ws.Range("E13").Value = 0.07
$ws.Range("F13").Value = 38.8

# Row 14
$ws.Range("A14").Value = "Heapq: heuristic = numMoves + cars between RedCar & exit"
$ws.Range("C14").Value = 0.88
$ws.Range("D14").Value = 0.1
$ws.Range("D14").NumberFormat = "0.00"
$ws.Range("E14").Value = 0.06
$ws.Range("F14").Value = 24.9

# Update selection to F11 as in the target sheet view
$ws.Range("F11").Select()
